$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update B2:B7 image paths to live under images/
$ws.Range("B2").Value = "images/up.jpg"
$ws.Range("B3").Value = "images/up.jpg"
$ws.Range("B4").Value = "images/down.jpg"
$ws.Range("B5").Value = "images/up.jpg"
$ws.Range("B6").Value = "images/down.jpg"
$ws.Range("B7").Value = "images/down.jpg"

# Unify the C column prompts into a single rating question
$ws.Range("C2").Value = "Rate your ability to control your brain"
$ws.Range("C3").Value = "Rate your ability to control your brain"
$ws.Range("C4").Value = "Rate your ability to control your brain"
$ws.Range("C5").Value = "Rate your ability to control your brain"
$ws.Range("C6").Value = "Rate your ability to control your brain"
$ws.Range("C7").Value = "Rate your ability to control your brain"

# Drop the now-unused opacity column (D)
$ws.Columns.Item(4).Delete()

# Move the active selection on Sheet1 from B14 to B8
$ws.Range("B8").Select()

# Nudge the workbook tab-split ratio (cosmetic window chrome)
$win = $wb.Windows.Item(1)
$win.TabRatio = 496

# Sheet2 / Sheet3: bump the page-layout zoom back to 100% and reset the
# printed "first page number" override (matches the fixed condition xmls)
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.PageSetup.Zoom = 100
$sheet2.PageSetup.FirstPageNumber = 0

$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet3.PageSetup.Zoom = 100
$sheet3.PageSetup.FirstPageNumber = 0
